{"js": "// Apply the \"Added many more features\" edit to the Great Rhino review.\n// Each entry is an exact, unique-enough search string paired with its\n// replacement. Title/heading text is searched with matchCase so we don't\n// accidentally hit unrelated runs.\nconst replacements = [\n  {\n    find: \"Play Great Rhino for Free: Review & Pros and Cons\",\n    replace: \"Play Great Rhino Free: Review of Features, Jackpots, and Visuals\",\n  },\n  {\n    find: \"Immersive African savannah graphics and design\",\n    replace: \"Intuitive and easy-to-play gameplay\",\n  },\n  {\n    find: \"Medium volatility for a balanced game experience\",\n    replace: \"Eye-catching and entertaining design\",\n  },\n  {\n    find: \"Standard RTP of 96.53%\",\n    replace: \"Stacked wilds and free spins feature\",\n  },\n  {\n    find: \"Two jackpots that can be won during free spins\",\n    replace: \"Chance to win significant rewards\",\n  },\n  {\n    find: \"No base game jackpot\",\n    replace: \"Medium volatility may not appeal to players seeking frequent wins\",\n  },\n  {\n    find: \"Free spins cannot be re-triggered\",\n    replace: \"Limited number of bonus features\",\n  },\n  {\n    find: \"Read our review of Great Rhino, a medium-volatility slot with 96.53% RTP. Enjoy two jackpots, stacked wilds and immersive graphics. Play for free!\",\n    replace: \"Explore the African savannah and play Great Rhino free. Find out about stacked wilds, free spins, and potential rewards in this review.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Great Rhino review.\n# Use Find/Replace over the whole document content for each exact phrase.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Great Rhino for Free: Review & Pros and Cons\"; Replace = \"Play Great Rhino Free: Review of Features, Jackpots, and Visuals\" },\n    @{ Find = \"Immersive African savannah graphics and design\"; Replace = \"Intuitive and easy-to-play gameplay\" },\n    @{ Find = \"Medium volatility for a balanced game experience\"; Replace = \"Eye-catching and entertaining design\" },\n    @{ Find = \"Standard RTP of 96.53%\"; Replace = \"Stacked wilds and free spins feature\" },\n    @{ Find = \"Two jackpots that can be won during free spins\"; Replace = \"Chance to win significant rewards\" },\n    @{ Find = \"No base game jackpot\"; Replace = \"Medium volatility may not appeal to players seeking frequent wins\" },\n    @{ Find = \"Free spins cannot be re-triggered\"; Replace = \"Limited number of bonus features\" },\n    @{ Find = \"Read our review of Great Rhino, a medium-volatility slot with 96.53% RTP. Enjoy two jackpots, stacked wilds and immersive graphics. Play for free!\"; Replace = \"Explore the African savannah and play Great Rhino free. Find out about stacked wilds, free spins, and potential rewards in this review.\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $r.Find\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Replacement.Text = $r.Replace\n    $range.Find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
